$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" quarterly sheet.
#    Clone the existing "2022-Q2" sheet (same column layout/formatting) and
#    drop the copy in front of it, then rename it and fill in the fresh
#    Q3 numbers.
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$existingQ2.Copy($existingQ2)
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# The D:G columns hold numeric-looking values that are stored as text in
# this workbook, so force text formatting before typing them in (otherwise
# Excel would auto-convert them to numbers).
$newQ3.Range("D2:G3").NumberFormat = "@"

$newQ3.Range("D2").Value = "0.13"
$newQ3.Range("E2").Value = "92.85"
$newQ3.Range("F2").Value = "3.29"
$newQ3.Range("G2").Value = "0.0043"
$newQ3.Range("H2").Value = 8

$newQ3.Range("D3").Value = "0.08"
$newQ3.Range("E3").Value = "92.85"
$newQ3.Range("F3").Value = "3.29"
$newQ3.Range("G3").Value = "0.0026"
$newQ3.Range("H3").Value = 8

# Drop the temporary text number-format now that the literal text values are
# committed, so the cells fall back to the same (unstyled) look the rest of
# the sheet uses.
$newQ3.Range("D2:G3").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) roll-up sheet: shift the quarter labels down
#    and append the new trailing 2022-Q1 row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("B3").Value = "2022-Q2"

$summary.Range("A3:D3").Copy()
$summary.Range("A4:D4").PasteSpecial(-4122) # xlPasteFormats
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.01

# Restore the originally-selected tab ("2022-Q1" stays the last sheet).
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Output "2022-Q3 sheet added and 总计 updated"
